# Updates the cryptos price/volume table (GitHub Actions scheduled refresh).
# Price values in column D are stored as plain text in the workbook (they
# use a "." thousands separator, e.g. "30.450.88", so Excel can't treat
# them as real numbers anyway). When the new price text IS a value Excel's
# parser would recognise as a number (e.g. "21.89"), a bare .Value
# assignment gets auto-converted to a float - so those writes are done with
# a leading "'" (forces text entry, exactly like typing '21.89 into the
# cell) and then ClearFormats() immediately strips the quote-prefix style
# Excel applies for that, leaving a plain text cell with no style override,
# matching the original formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "30.450.88"
$ws.Cells.Item(2,5).Value = "  +0.96%  "
$ws.Cells.Item(3,4).Value = "1.878.01"
$ws.Cells.Item(3,5).Value = "  +0.80%  "
$ws.Cells.Item(5,4).Value = "'247.27"
$ws.Cells.Item(5,4).ClearFormats() | Out-Null
$ws.Cells.Item(5,5).Value = "  +5.62%  "
$ws.Cells.Item(6,5).Value = "  -0.08%  "
$ws.Cells.Item(7,4).Value = "'0.4762"
$ws.Cells.Item(7,4).ClearFormats() | Out-Null
$ws.Cells.Item(7,5).Value = "  +1.65%  "
$ws.Cells.Item(8,4).Value = "'0.2899"
$ws.Cells.Item(8,4).ClearFormats() | Out-Null
$ws.Cells.Item(8,5).Value = "  +1.42%  "
$ws.Cells.Item(9,4).Value = "'0.06530"
$ws.Cells.Item(9,4).ClearFormats() | Out-Null
$ws.Cells.Item(9,5).Value = "  +0.60%  "
$ws.Cells.Item(10,4).Value = "'21.89"
$ws.Cells.Item(10,4).ClearFormats() | Out-Null
$ws.Cells.Item(10,5).Value = "  +3.62%  "
$ws.Cells.Item(11,4).Value = "'0.07726"
$ws.Cells.Item(11,4).ClearFormats() | Out-Null
$ws.Cells.Item(11,5).Value = "  -0.52%  "
$ws.Cells.Item(12,4).Value = "'97.02"
$ws.Cells.Item(12,4).ClearFormats() | Out-Null
$ws.Cells.Item(12,5).Value = "  +3.43%  "
$ws.Cells.Item(13,4).Value = "'0.7375"
$ws.Cells.Item(13,4).ClearFormats() | Out-Null
$ws.Cells.Item(13,5).Value = "  +7.91%  "
$ws.Cells.Item(14,4).Value = "1.879.86"
$ws.Cells.Item(14,5).Value = "  +1.06%  "
$ws.Cells.Item(15,4).Value = "'5.131"
$ws.Cells.Item(15,4).ClearFormats() | Out-Null
$ws.Cells.Item(15,5).Value = "  +1.56%  "
$ws.Cells.Item(16,4).Value = "'272.36"
$ws.Cells.Item(16,4).ClearFormats() | Out-Null
$ws.Cells.Item(16,5).Value = "  +1.14%  "
$ws.Cells.Item(17,4).Value = "30.450.35"
$ws.Cells.Item(17,5).Value = "  +0.99%  "
$ws.Cells.Item(18,5).Value = "  +2.39%  "
$ws.Cells.Item(19,4).Value = "'0.000007609"
$ws.Cells.Item(19,4).ClearFormats() | Out-Null
$ws.Cells.Item(20,5).Value = "  -0.07%  "
$ws.Cells.Item(21,4).Value = "2.125.34"
$ws.Cells.Item(21,5).Value = "  +1.08%  "
$ws.Cells.Item(22,5).Value = "  +0.00%  "
$ws.Cells.Item(23,4).Value = "'5.254"
$ws.Cells.Item(23,4).ClearFormats() | Out-Null
$ws.Cells.Item(23,5).Value = "  +2.07%  "
$ws.Cells.Item(24,4).Value = "'6.188"
$ws.Cells.Item(24,4).ClearFormats() | Out-Null
$ws.Cells.Item(24,5).Value = "  +1.37%  "
$ws.Cells.Item(25,4).Value = "'9.332"
$ws.Cells.Item(25,4).ClearFormats() | Out-Null
$ws.Cells.Item(25,5).Value = "  -0.08%  "
$ws.Cells.Item(26,4).Value = "'164.09"
$ws.Cells.Item(26,4).ClearFormats() | Out-Null
$ws.Cells.Item(26,5).Value = "  -0.83%  "
$ws.Cells.Item(27,4).Value = "'18.88"
$ws.Cells.Item(27,4).ClearFormats() | Out-Null
$ws.Cells.Item(27,5).Value = "  +1.79%  "
$ws.Cells.Item(28,4).Value = "'1.942"
$ws.Cells.Item(28,4).ClearFormats() | Out-Null
$ws.Cells.Item(28,5).Value = "  +2.70%  "
$ws.Cells.Item(29,4).Value = "'1.372"
$ws.Cells.Item(29,4).ClearFormats() | Out-Null
$ws.Cells.Item(29,5).Value = "  +0.66%  "
$ws.Cells.Item(30,4).Value = "'0.09949"
$ws.Cells.Item(30,4).ClearFormats() | Out-Null
$ws.Cells.Item(30,5).Value = "  -0.03%  "
$ws.Cells.Item(31,4).Value = "'1.520"
$ws.Cells.Item(31,4).ClearFormats() | Out-Null
$ws.Cells.Item(31,5).Value = "  +4.80%  "
$ws.Cells.Item(32,4).Value = "'4.314"
$ws.Cells.Item(32,4).ClearFormats() | Out-Null
$ws.Cells.Item(32,5).Value = "  +2.11%  "
$ws.Cells.Item(33,4).Value = "'4.071"
$ws.Cells.Item(33,4).ClearFormats() | Out-Null
$ws.Cells.Item(34,4).Value = "'0.04777"
$ws.Cells.Item(34,4).ClearFormats() | Out-Null
$ws.Cells.Item(34,5).Value = "  +2.04%  "
$ws.Cells.Item(35,4).Value = "'1.126"
$ws.Cells.Item(35,4).ClearFormats() | Out-Null
$ws.Cells.Item(35,5).Value = "  +0.74%  "
$ws.Cells.Item(36,4).Value = "'0.7025"
$ws.Cells.Item(36,4).ClearFormats() | Out-Null
$ws.Cells.Item(36,5).Value = "  +1.74%  "
$ws.Cells.Item(37,4).Value = "'2.714"
$ws.Cells.Item(37,4).ClearFormats() | Out-Null
$ws.Cells.Item(37,5).Value = "  +0.44%  "
$ws.Cells.Item(38,4).Value = "'0.01867"
$ws.Cells.Item(38,4).ClearFormats() | Out-Null
$ws.Cells.Item(38,5).Value = "  +1.61%  "
$ws.Cells.Item(39,4).Value = "'2.730"
$ws.Cells.Item(39,4).ClearFormats() | Out-Null
$ws.Cells.Item(39,5).Value = "  -1.02%  "
$ws.Cells.Item(40,4).Value = "'6.344"
$ws.Cells.Item(40,4).ClearFormats() | Out-Null
$ws.Cells.Item(40,5).Value = "  +0.10%  "
$ws.Cells.Item(41,2).Value = "RenderToken"
$ws.Cells.Item(41,3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(41,4).Value = "'1.949"
$ws.Cells.Item(41,4).ClearFormats() | Out-Null
$ws.Cells.Item(41,5).Value = "  +2.92%  "
$ws.Cells.Item(42,2).Value = "Aave"
$ws.Cells.Item(42,3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(42,4).Value = "'70.74"
$ws.Cells.Item(42,4).ClearFormats() | Out-Null
$ws.Cells.Item(42,5).Value = "  -0.92%  "
$ws.Cells.Item(43,4).Value = "'0.4196"
$ws.Cells.Item(43,4).ClearFormats() | Out-Null
$ws.Cells.Item(43,5).Value = "  +3.45%  "
$ws.Cells.Item(44,4).Value = "'0.9998"
$ws.Cells.Item(44,4).ClearFormats() | Out-Null
$ws.Cells.Item(45,4).Value = "'0.8371"
$ws.Cells.Item(45,4).ClearFormats() | Out-Null
$ws.Cells.Item(45,5).Value = "  +0.30%  "
$ws.Cells.Item(46,4).Value = "'102.91"
$ws.Cells.Item(46,4).ClearFormats() | Out-Null
$ws.Cells.Item(46,5).Value = "  +0.71%  "
$ws.Cells.Item(47,4).Value = "'9.262"
$ws.Cells.Item(47,4).ClearFormats() | Out-Null
$ws.Cells.Item(47,5).Value = "  +1.45%  "
$ws.Cells.Item(48,4).Value = "'7.095"
$ws.Cells.Item(48,4).ClearFormats() | Out-Null
$ws.Cells.Item(48,5).Value = "  +2.01%  "
$ws.Cells.Item(49,4).Value = "'35.56"
$ws.Cells.Item(49,4).ClearFormats() | Out-Null
$ws.Cells.Item(49,5).Value = "  +4.28%  "
$ws.Cells.Item(50,4).Value = "'928.36"
$ws.Cells.Item(50,4).ClearFormats() | Out-Null
$ws.Cells.Item(50,5).Value = "  -0.50%  "
$ws.Cells.Item(51,4).Value = "'0.05642"
$ws.Cells.Item(51,4).ClearFormats() | Out-Null
$ws.Cells.Item(51,5).Value = "  +1.16%  "
